$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new row 6 with the latest scraped tracker result -------------

$ws.Range("A6").Value = 14494999

# "2025-09-06" looks like a date to Excel's auto-detection; build it as a
# formula-derived literal string and paste the computed value back so the
# cell keeps a genuine text value (and not a date serial) without leaving
# a helper number-format style behind.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""2025-09-06"""
$scratch.Copy()
$ws.Range("B6").PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("C6").Value = "Aryna Sabalenka"
$ws.Range("D6").Value = "Amanda Anisimova"
$ws.Range("E6").Value = "Gana Aryna Sabalenka"
$ws.Range("F6").Value = 1.48

# resultado / profit are still blank for this just-added fixture; touch the
# cells so the row materializes through column H like the source tracker.
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Style = "Normal"
